# docpac_oct15 edit
#   - Remove Duplicate Bullet Point
#   - Removed "ExpressJS POST Requests"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the duplicate "ExpressJS POST Requests " bullet paragraph
#    that follows "Notebook 6: ExpressJS POST Requests". Selecting the
#    whole paragraph (including its end-of-paragraph mark) and deleting
#    it mirrors what a user would do; Word then drops a _GoBack bookmark
#    at the edit point.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ExpressJS POST Requests `r") {
        $editSpot = $p.Range.Duplicate
        $editSpot.Collapse(1)
        $p.Range.Delete()
        $d.Bookmarks.Add("_GoBack", $editSpot)
        break
    }
}

# ---------------------------------------------------------------------
# 2) A handful of paragraphs had their text split across two runs with
#    identical formatting. Running them back through Find & Replace
#    (old text -> same text) collapses each pair back into a single run,
#    matching the cleaned-up document.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Required Documentati" + "on", $false, $false, $false, $false, $false, $true, 1, $false, "Required Documentation", 2) | Out-Null

$d.Content.Find.Execute("Write a new Notebook entry explaining to a new programmer how to submit form data in an HTML file in your web browser. Must contain the following:", $false, $false, $false, $false, $false, $true, 1, $false, "Write a new Notebook entry explaining to a new programmer how to submit form data in an HTML file in your web browser. Must contain the following:", 2) | Out-Null

$d.Content.Find.Execute("Right now, your team website has the “/feedback” endpoint. This can collect Query Parameters in the URL to save data to a JSON. There are many drawbacks to this. You will need to upgrade your team website to send a user to a new web page when they send a GET request with no Query Parameters to the endpoint “/feedback” This web page will have:", $false, $false, $false, $false, $false, $true, 1, $false, "Right now, your team website has the “/feedback” endpoint. This can collect Query Parameters in the URL to save data to a JSON. There are many drawbacks to this. You will need to upgrade your team website to send a user to a new web page when they send a GET request with no Query Parameters to the endpoint “/feedback” This web page will have:", 2) | Out-Null

$d.Content.Find.Execute("When the user clicks submit, your server will handle the “/feedback” POST request and read the form data. If it is valid (contains a username and comment), it will push to the comments array and save it to the comments file the same way your “/feedback” GET request does now.", $false, $false, $false, $false, $false, $true, 1, $false, "When the user clicks submit, your server will handle the “/feedback” POST request and read the form data. If it is valid (contains a username and comment), it will push to the comments array and save it to the comments file the same way your “/feedback” GET request does now.", 2) | Out-Null
